$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.219.33"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "'1.904.25"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'307.90"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.5263"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").Value = "'0.3827"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("D9").Value = "'0.07304"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'21.57"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D12").Value = "'0.08080"
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("D13").Value = "'96.02"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'5.369"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'1.765.63"
$ws.Range("E15").Value = "  -7.21%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'0.000008681"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'14.74"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'27.255.66"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "'2.338"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "'150.20"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "'116.81"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "'4.850"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "'4.879"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'0.09226"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'0.8191"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").Value = "'0.05077"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "'1.231"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "'2.984"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").Value = "'3.359"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("D37").Value = "'2.723"
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("D38").Value = "'0.5736"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").Value = "'0.01997"
$ws.Range("D41").Value = "'9.030"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'6.627"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'117.05"
$ws.Range("D44").Value = "'0.1525"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "'0.4933"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").Value = "'1.003"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "'1.641"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "'38.63"
$ws.Range("E49").Value = "  +3.15%  "
$ws.Range("D50").Value = "'64.08"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "'0.05969"
$ws.Range("E51").Value = "  +0.26%  "
